# Update metrics table (columns B:Q, rows 2:26) with new values
# produced by retraining the model ("atualizado todo o treinamento para o novo lm").
# Note: values are written in plain decimal (not scientific notation) because
# the COM-interop script parser does not accept exponent literals like 1e-07.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.9999989698596451,
    0.9990244665669379,
    0.9999987405318091,
    0.9999999999998166,
    0.9999995167451518,
    0.0000009615908033740045,
    0.0009106176387947784,
    0.000001066751762353674,
    0.0000000000002489308586522353,
    0.0000005333760056422664,
    0.00005130350677945762,
    0.000980607364531801,
    1.000024723368518,
    0.001022353850177004,
    77.7093536741803,
    108.1812492958853
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
